# Updates the cryptos price/volume table to the refreshed values from the
# GitHub Actions data pull. Column D ("Price") values are stored as text
# (e.g. "30.241.90", "0.4670", "1.000") -- many aren't even valid numbers
# (multiple '.' used as thousands separators) and the ones that are would
# lose meaningful trailing zeros / formatting if Excel auto-converted them
# to numeric values. Prefixing with a leading apostrophe (the same trick
# used when typing a number-looking value into Excel by hand) forces the
# cell to keep its text type, matching the original inline/shared string
# representation exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''30.241.90'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '''1.861.85'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''235.71'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').Value = '''0.4670'
$ws.Range('E7').Value = '  -0.58%  '
$ws.Range('D8').Value = '''0.2829'
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('D9').Value = '''0.06510'
$ws.Range('E9').Value = '  -0.60%  '
$ws.Range('D10').Value = '''21.39'
$ws.Range('E10').Value = '  +6.33%  '
$ws.Range('D12').Value = '''96.96'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').Value = '''1.864.22'
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').Value = '''5.145'
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('D15').Value = '''0.6770'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').Value = '''278.01'
$ws.Range('E16').Value = '  -2.11%  '
$ws.Range('D17').Value = '''30.242.77'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = '''13.69'
$ws.Range('E18').Value = '  +8.89%  '
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').Value = '''5.380'
$ws.Range('E20').Value = '  -1.15%  '
$ws.Range('D21').Value = '''2.108.88'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').Value = '''0.000007305'
$ws.Range('E22').Value = '  +0.93%  '
$ws.Range('D23').Value = '''1.001'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '''167.30'
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('D26').Value = '''9.134'
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('E29').Value = '  +3.13%  '
$ws.Range('D30').Value = '''0.09699'
$ws.Range('E30').Value = '  +1.09%  '
$ws.Range('D31').Value = '''4.363'
$ws.Range('E31').Value = '  -1.21%  '
$ws.Range('D32').Value = '''1.478'
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('D33').Value = '''4.017'
$ws.Range('E33').Value = '  -1.99%  '
$ws.Range('D34').Value = '''0.04719'
$ws.Range('E34').Value = '  +1.00%  '
$ws.Range('E35').Value = '  +2.07%  '
$ws.Range('D36').Value = '''0.7036'
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').Value = '''0.01856'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = '''2.580'
$ws.Range('E39').Value = '  +2.77%  '
$ws.Range('D40').Value = '''6.313'
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('D41').Value = '''75.16'
$ws.Range('E41').Value = '  +4.32%  '
$ws.Range('D42').Value = '''1.954'
$ws.Range('E42').Value = '  +0.88%  '
$ws.Range('D43').Value = '''0.8486'
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '''1.000'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '''0.4161'
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').Value = '''103.35'
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range('D47').Value = '''979.69'
$ws.Range('E47').Value = '  -2.49%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '''7.143'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''9.268'
$ws.Range('E49').Value = '  +2.74%  '
$ws.Range('D50').Value = '''33.99'
$ws.Range('E50').Value = '  +0.82%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.05649'
$ws.Range('E51').Value = '  +0.09%  '
